$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# The AddCustomerTest data block (rows 1-4) stays as-is. The repeated
# duplicate rows that used to follow it (old rows 5-13) are removed, and
# the OpenAccountTest block (old rows 14-17) is pulled up to rows 6-9,
# leaving row 5 blank as a separator - matching the trimmed-down test data.
$ws.Range("A5:E17").ClearContents()

$ws.Range("A6").Value = "OpenAccountTest"

$ws.Range("A7").Value = "Runmode"
$ws.Range("B7").Value = "customer"
$ws.Range("C7").Value = "currency"

$ws.Range("A8").Value = "Y"
$ws.Range("B8").Value = "manish k"
$ws.Range("C8").Value = "Rupee"

$ws.Range("A9").Value = "N"
$ws.Range("B9").Value = "jyoti k"
$ws.Range("C9").Value = "Dollar"
